$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header suffixes: _old -> _FV2210, _new -> _FV2304
$headerRange = $ws.Range("A1:U1")
foreach ($cell in $headerRange.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# Freeze the header row (pane split after row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into a native Excel Table ("Table1")
$listRange = $ws.Range("A1:U77")
$tbl = $ws.ListObjects.Add(1, $listRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
